# This commit is a pure OOXML attribute-canonicalization / serialization-order
# change (the diff reorders XML attributes alphabetically throughout
# ppt/slideMasters/slideMaster1.xml -- e.g. <a:bodyPr vert="horz" lIns=...>
# becomes <a:bodyPr anchor="ctr" bIns=...> -- and reorders the theme's
# position when the archive is rebuilt). It carries no semantic/visual
# change: no text, geometry, formatting, or structural content differs
# between the two sides of the diff.
#
# There is no PowerPoint object-model operation that edits raw XML
# attribute ordering, and the slide master already contains exactly the
# content described by both sides of the diff, so no edit is necessary
# here -- we simply touch the presentation object model to confirm the
# document is accessible, without altering any content.
$p = $ppt.ActivePresentation
$null = $p.SlideMaster
